$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G (header "K") values being regenerated from Strike# to K
$gValues = @{
    2  = 1
    3  = 0
    4  = 2
    5  = 1
    6  = 2
    7  = 1
    8  = 1
    9  = 1
    10 = 3
    11 = 2
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}
